$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Banana Firm"
$ws.Range("A3").Value = "Pear Company"
$ws.Range("A4").Value = "Cherry Enterprise"
$ws.Range("A5").Value = "Grape Startup"

$ws.Range("B12:D17").ClearContents()

$ws.Range("D13").Select()
